# Apply updated cryptocurrency price/volume figures to columns D and E.
# Values that look like plain numbers must be forced to Text format first
# so Excel keeps them as literal strings (matching the source data, which
# stores prices/volumes as text, not numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.333.43"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.940.84"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "494.40"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.09"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.731"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000350"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.26"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "4.572.32"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "3.943.16"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.25"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.90"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "69.344.62"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.14"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.56"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.04"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.10"
$ws.Range("E25").Value = "  +9.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.79"
$ws.Range("E26").Value = "  +3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.21"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "702.57"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.39"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.464"
$ws.Range("E34").Value = "  +15.87%  "
$ws.Range("D35").Value = "0.0₃0893"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "61.91"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "40.86"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.150"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0489"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  +7.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.99"
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.68"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "0.0₆0340"
$ws.Range("E51").Value = "  -5.39%  "
